$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-09-18 Thursday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-09-19 Friday", 2) | Out-Null

# Update the table of arithmetic problems (20 rows x 5 columns)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Find.Execute("73-24=", $false, $false, $false, $false, $false, $true, 1, $false, "84-61=", 2) | Out-Null
$t.Cell(1, 2).Range.Find.Execute("80-47=", $false, $false, $false, $false, $false, $true, 1, $false, "25+6=", 2) | Out-Null
$t.Cell(1, 3).Range.Find.Execute("71-24=", $false, $false, $false, $false, $false, $true, 1, $false, "31-22=", 2) | Out-Null
$t.Cell(1, 4).Range.Find.Execute("71-16=", $false, $false, $false, $false, $false, $true, 1, $false, "30+69=", 2) | Out-Null
$t.Cell(1, 5).Range.Find.Execute("45+39=", $false, $false, $false, $false, $false, $true, 1, $false, "79-8=", 2) | Out-Null

$t.Cell(2, 1).Range.Find.Execute("26+35=", $false, $false, $false, $false, $false, $true, 1, $false, "46+19=", 2) | Out-Null
$t.Cell(2, 2).Range.Find.Execute("6+35=", $false, $false, $false, $false, $false, $true, 1, $false, "92-0=", 2) | Out-Null
$t.Cell(2, 3).Range.Find.Execute("71-42=", $false, $false, $false, $false, $false, $true, 1, $false, "88-35=", 2) | Out-Null
$t.Cell(2, 4).Range.Find.Execute("80-10=", $false, $false, $false, $false, $false, $true, 1, $false, "11+41=", 2) | Out-Null
$t.Cell(2, 5).Range.Find.Execute("35+33=", $false, $false, $false, $false, $false, $true, 1, $false, "58+16=", 2) | Out-Null

$t.Cell(3, 1).Range.Find.Execute("27+34=", $false, $false, $false, $false, $false, $true, 1, $false, "75-7=", 2) | Out-Null
$t.Cell(3, 2).Range.Find.Execute("75-17=", $false, $false, $false, $false, $false, $true, 1, $false, "44-2=", 2) | Out-Null
$t.Cell(3, 3).Range.Find.Execute("33+57=", $false, $false, $false, $false, $false, $true, 1, $false, "17+34=", 2) | Out-Null
$t.Cell(3, 4).Range.Find.Execute("36-34=", $false, $false, $false, $false, $false, $true, 1, $false, "77-42=", 2) | Out-Null
$t.Cell(3, 5).Range.Find.Execute("10+68=", $false, $false, $false, $false, $false, $true, 1, $false, "95-77=", 2) | Out-Null

$t.Cell(4, 1).Range.Find.Execute("80-68=", $false, $false, $false, $false, $false, $true, 1, $false, "94-38=", 2) | Out-Null
$t.Cell(4, 2).Range.Find.Execute("66-41=", $false, $false, $false, $false, $false, $true, 1, $false, "30-5=", 2) | Out-Null
$t.Cell(4, 3).Range.Find.Execute("84-67=", $false, $false, $false, $false, $false, $true, 1, $false, "86-66=", 2) | Out-Null
$t.Cell(4, 4).Range.Find.Execute("94-28=", $false, $false, $false, $false, $false, $true, 1, $false, "75-68=", 2) | Out-Null
$t.Cell(4, 5).Range.Find.Execute("35+15=", $false, $false, $false, $false, $false, $true, 1, $false, "6+8=", 2) | Out-Null

$t.Cell(5, 1).Range.Find.Execute("89-33=", $false, $false, $false, $false, $false, $true, 1, $false, "70-10=", 2) | Out-Null
$t.Cell(5, 2).Range.Find.Execute("99-68=", $false, $false, $false, $false, $false, $true, 1, $false, "48+41=", 2) | Out-Null
$t.Cell(5, 3).Range.Find.Execute("34+10=", $false, $false, $false, $false, $false, $true, 1, $false, "95-32=", 2) | Out-Null
$t.Cell(5, 4).Range.Find.Execute("22+1=", $false, $false, $false, $false, $false, $true, 1, $false, "23-6=", 2) | Out-Null
$t.Cell(5, 5).Range.Find.Execute("30-22=", $false, $false, $false, $false, $false, $true, 1, $false, "35+47=", 2) | Out-Null

$t.Cell(6, 1).Range.Find.Execute("10-2=", $false, $false, $false, $false, $false, $true, 1, $false, "78-76=", 2) | Out-Null
$t.Cell(6, 2).Range.Find.Execute("17+53=", $false, $false, $false, $false, $false, $true, 1, $false, "53+30=", 2) | Out-Null
$t.Cell(6, 3).Range.Find.Execute("21+57=", $false, $false, $false, $false, $false, $true, 1, $false, "12+25=", 2) | Out-Null
$t.Cell(6, 4).Range.Find.Execute("10+36=", $false, $false, $false, $false, $false, $true, 1, $false, "25+62=", 2) | Out-Null
$t.Cell(6, 5).Range.Find.Execute("48-12=", $false, $false, $false, $false, $false, $true, 1, $false, "61-33=", 2) | Out-Null

$t.Cell(7, 1).Range.Find.Execute("17-15=", $false, $false, $false, $false, $false, $true, 1, $false, "14+37=", 2) | Out-Null
$t.Cell(7, 2).Range.Find.Execute("34+5=", $false, $false, $false, $false, $false, $true, 1, $false, "39+1=", 2) | Out-Null
$t.Cell(7, 3).Range.Find.Execute("77+22=", $false, $false, $false, $false, $false, $true, 1, $false, "96-70=", 2) | Out-Null
$t.Cell(7, 4).Range.Find.Execute("82-8=", $false, $false, $false, $false, $false, $true, 1, $false, "73-12=", 2) | Out-Null
$t.Cell(7, 5).Range.Find.Execute("22+11=", $false, $false, $false, $false, $false, $true, 1, $false, "94-0=", 2) | Out-Null

$t.Cell(8, 1).Range.Find.Execute("4+5=", $false, $false, $false, $false, $false, $true, 1, $false, "92-53=", 2) | Out-Null
$t.Cell(8, 2).Range.Find.Execute("90+9=", $false, $false, $false, $false, $false, $true, 1, $false, "75-6=", 2) | Out-Null
$t.Cell(8, 3).Range.Find.Execute("64-24=", $false, $false, $false, $false, $false, $true, 1, $false, "69+27=", 2) | Out-Null
$t.Cell(8, 4).Range.Find.Execute("7+65=", $false, $false, $false, $false, $false, $true, 1, $false, "96-48=", 2) | Out-Null
$t.Cell(8, 5).Range.Find.Execute("23+11=", $false, $false, $false, $false, $false, $true, 1, $false, "33+58=", 2) | Out-Null

$t.Cell(9, 1).Range.Find.Execute("4+52=", $false, $false, $false, $false, $false, $true, 1, $false, "16+26=", 2) | Out-Null
$t.Cell(9, 2).Range.Find.Execute("98-69=", $false, $false, $false, $false, $false, $true, 1, $false, "8+51=", 2) | Out-Null
$t.Cell(9, 3).Range.Find.Execute("84-0=", $false, $false, $false, $false, $false, $true, 1, $false, "41+41=", 2) | Out-Null
$t.Cell(9, 4).Range.Find.Execute("67-66=", $false, $false, $false, $false, $false, $true, 1, $false, "39+23=", 2) | Out-Null
$t.Cell(9, 5).Range.Find.Execute("48-46=", $false, $false, $false, $false, $false, $true, 1, $false, "24+35=", 2) | Out-Null

$t.Cell(10, 1).Range.Find.Execute("79-7=", $false, $false, $false, $false, $false, $true, 1, $false, "99-59=", 2) | Out-Null
$t.Cell(10, 2).Range.Find.Execute("44+25=", $false, $false, $false, $false, $false, $true, 1, $false, "88+10=", 2) | Out-Null
$t.Cell(10, 3).Range.Find.Execute("22+76=", $false, $false, $false, $false, $false, $true, 1, $false, "19+10=", 2) | Out-Null
$t.Cell(10, 4).Range.Find.Execute("80-13=", $false, $false, $false, $false, $false, $true, 1, $false, "0+60=", 2) | Out-Null
$t.Cell(10, 5).Range.Find.Execute("96-10=", $false, $false, $false, $false, $false, $true, 1, $false, "46-36=", 2) | Out-Null

$t.Cell(11, 1).Range.Find.Execute("28+52=", $false, $false, $false, $false, $false, $true, 1, $false, "48+26=", 2) | Out-Null
$t.Cell(11, 2).Range.Find.Execute("70+22=", $false, $false, $false, $false, $false, $true, 1, $false, "81-79=", 2) | Out-Null
$t.Cell(11, 3).Range.Find.Execute("54-21=", $false, $false, $false, $false, $false, $true, 1, $false, "53-52=", 2) | Out-Null
$t.Cell(11, 4).Range.Find.Execute("0+61=", $false, $false, $false, $false, $false, $true, 1, $false, "30-25=", 2) | Out-Null
$t.Cell(11, 5).Range.Find.Execute("83-16=", $false, $false, $false, $false, $false, $true, 1, $false, "84-77=", 2) | Out-Null

$t.Cell(12, 1).Range.Find.Execute("81-58=", $false, $false, $false, $false, $false, $true, 1, $false, "2+94=", 2) | Out-Null
$t.Cell(12, 2).Range.Find.Execute("94-30=", $false, $false, $false, $false, $false, $true, 1, $false, "16+66=", 2) | Out-Null
$t.Cell(12, 3).Range.Find.Execute("36+15=", $false, $false, $false, $false, $false, $true, 1, $false, "1+70=", 2) | Out-Null
$t.Cell(12, 4).Range.Find.Execute("30+65=", $false, $false, $false, $false, $false, $true, 1, $false, "82-20=", 2) | Out-Null
$t.Cell(12, 5).Range.Find.Execute("58-3=", $false, $false, $false, $false, $false, $true, 1, $false, "14-4=", 2) | Out-Null

$t.Cell(13, 1).Range.Find.Execute("22+44=", $false, $false, $false, $false, $false, $true, 1, $false, "78-37=", 2) | Out-Null
$t.Cell(13, 2).Range.Find.Execute("15+13=", $false, $false, $false, $false, $false, $true, 1, $false, "3+83=", 2) | Out-Null
$t.Cell(13, 3).Range.Find.Execute("81-41=", $false, $false, $false, $false, $false, $true, 1, $false, "10+13=", 2) | Out-Null
$t.Cell(13, 4).Range.Find.Execute("57+23=", $false, $false, $false, $false, $false, $true, 1, $false, "10+43=", 2) | Out-Null
$t.Cell(13, 5).Range.Find.Execute("22-13=", $false, $false, $false, $false, $false, $true, 1, $false, "97-12=", 2) | Out-Null

$t.Cell(14, 1).Range.Find.Execute("80-79=", $false, $false, $false, $false, $false, $true, 1, $false, "74+23=", 2) | Out-Null
$t.Cell(14, 2).Range.Find.Execute("62-46=", $false, $false, $false, $false, $false, $true, 1, $false, "30+46=", 2) | Out-Null
$t.Cell(14, 3).Range.Find.Execute("3+41=", $false, $false, $false, $false, $false, $true, 1, $false, "51-31=", 2) | Out-Null
$t.Cell(14, 4).Range.Find.Execute("82+14=", $false, $false, $false, $false, $false, $true, 1, $false, "41+58=", 2) | Out-Null
$t.Cell(14, 5).Range.Find.Execute("20+12=", $false, $false, $false, $false, $false, $true, 1, $false, "17+64=", 2) | Out-Null

$t.Cell(15, 1).Range.Find.Execute("99-12=", $false, $false, $false, $false, $false, $true, 1, $false, "64-61=", 2) | Out-Null
$t.Cell(15, 2).Range.Find.Execute("11+43=", $false, $false, $false, $false, $false, $true, 1, $false, "57+10=", 2) | Out-Null
$t.Cell(15, 3).Range.Find.Execute("12+61=", $false, $false, $false, $false, $false, $true, 1, $false, "63+27=", 2) | Out-Null
$t.Cell(15, 4).Range.Find.Execute("41+51=", $false, $false, $false, $false, $false, $true, 1, $false, "48+40=", 2) | Out-Null
$t.Cell(15, 5).Range.Find.Execute("90-38=", $false, $false, $false, $false, $false, $true, 1, $false, "70+4=", 2) | Out-Null

$t.Cell(16, 1).Range.Find.Execute("40-36=", $false, $false, $false, $false, $false, $true, 1, $false, "6+63=", 2) | Out-Null
$t.Cell(16, 2).Range.Find.Execute("23+47=", $false, $false, $false, $false, $false, $true, 1, $false, "60+1=", 2) | Out-Null
$t.Cell(16, 3).Range.Find.Execute("10+52=", $false, $false, $false, $false, $false, $true, 1, $false, "5+82=", 2) | Out-Null
$t.Cell(16, 4).Range.Find.Execute("53+12=", $false, $false, $false, $false, $false, $true, 1, $false, "26-18=", 2) | Out-Null
$t.Cell(16, 5).Range.Find.Execute("27-26=", $false, $false, $false, $false, $false, $true, 1, $false, "75-30=", 2) | Out-Null

$t.Cell(17, 1).Range.Find.Execute("37-2=", $false, $false, $false, $false, $false, $true, 1, $false, "66-32=", 2) | Out-Null
$t.Cell(17, 2).Range.Find.Execute("76-65=", $false, $false, $false, $false, $false, $true, 1, $false, "79-38=", 2) | Out-Null
$t.Cell(17, 3).Range.Find.Execute("6+1=", $false, $false, $false, $false, $false, $true, 1, $false, "58+16=", 2) | Out-Null
$t.Cell(17, 4).Range.Find.Execute("49+34=", $false, $false, $false, $false, $false, $true, 1, $false, "10+6=", 2) | Out-Null
$t.Cell(17, 5).Range.Find.Execute("96-62=", $false, $false, $false, $false, $false, $true, 1, $false, "25+61=", 2) | Out-Null

$t.Cell(18, 1).Range.Find.Execute("31+12=", $false, $false, $false, $false, $false, $true, 1, $false, "31+16=", 2) | Out-Null
$t.Cell(18, 2).Range.Find.Execute("12+48=", $false, $false, $false, $false, $false, $true, 1, $false, "61-55=", 2) | Out-Null
$t.Cell(18, 3).Range.Find.Execute("81-66=", $false, $false, $false, $false, $false, $true, 1, $false, "78-25=", 2) | Out-Null
$t.Cell(18, 4).Range.Find.Execute("53-6=", $false, $false, $false, $false, $false, $true, 1, $false, "5+42=", 2) | Out-Null
$t.Cell(18, 5).Range.Find.Execute("42+47=", $false, $false, $false, $false, $false, $true, 1, $false, "82+5=", 2) | Out-Null

$t.Cell(19, 1).Range.Find.Execute("76-39=", $false, $false, $false, $false, $false, $true, 1, $false, "41-33=", 2) | Out-Null
$t.Cell(19, 2).Range.Find.Execute("38+0=", $false, $false, $false, $false, $false, $true, 1, $false, "61-31=", 2) | Out-Null
$t.Cell(19, 3).Range.Find.Execute("97-79=", $false, $false, $false, $false, $false, $true, 1, $false, "2+39=", 2) | Out-Null
$t.Cell(19, 4).Range.Find.Execute("87+6=", $false, $false, $false, $false, $false, $true, 1, $false, "66+11=", 2) | Out-Null
$t.Cell(19, 5).Range.Find.Execute("18-8=", $false, $false, $false, $false, $false, $true, 1, $false, "23+40=", 2) | Out-Null

$t.Cell(20, 1).Range.Find.Execute("3+43=", $false, $false, $false, $false, $false, $true, 1, $false, "73-15=", 2) | Out-Null
$t.Cell(20, 2).Range.Find.Execute("83+10=", $false, $false, $false, $false, $false, $true, 1, $false, "10+24=", 2) | Out-Null
$t.Cell(20, 3).Range.Find.Execute("19+4=", $false, $false, $false, $false, $false, $true, 1, $false, "39-17=", 2) | Out-Null
$t.Cell(20, 4).Range.Find.Execute("19+6=", $false, $false, $false, $false, $false, $true, 1, $false, "25+59=", 2) | Out-Null
$t.Cell(20, 5).Range.Find.Execute("15-11=", $false, $false, $false, $false, $false, $true, 1, $false, "72-13=", 2) | Out-Null

